# Edit: after the "$ git pull origin master" paragraph there is an empty
# paragraph followed by a paragraph that only holds the "_GoBack" bookmark.
# The commit removes the empty paragraph and instead prefixes the bookmark
# paragraph with the runs "$ " + "git" (spell-check wrapped) + " ".

$d = $word.ActiveDocument

# Locate the paragraph that contains the "pull origin master" command text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pull origin master*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not locate the '...pull origin master' paragraph"
}

# The paragraph immediately following it should be empty - this is the one
# the diff removes.
$emptyPara = $target.Next()
if ($emptyPara.Range.Text.Trim().Length -ne 0) {
    throw "Expected an empty paragraph after 'pull origin master', got: [$($emptyPara.Range.Text)]"
}

# Deleting its range merges the empty paragraph away, leaving the bookmark
# paragraph directly after $target.
$emptyPara.Range.Delete()

# Re-fetch the (now adjacent) paragraph that holds the _GoBack bookmark -
# object references obtained before the Delete() call can be stale, so look
# it up again via Next().
$bookmarkPara = $target.Next()

# Collapsed range at the very start of that paragraph, i.e. right before the
# bookmarkStart element.
$insertPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = '<w:p ' + $wNs + '>' +
       '<w:r><w:t xml:space="preserve">$ </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>git</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '</w:p>'

# InsertXML on a range wrapped in <w:p> merges the runs into the existing
# paragraph at the insertion point instead of creating a new paragraph.
$insertPoint.InsertXML($xml)
